# Portfolio/Tables.xlsx — "added broken continue button"
# Adds two new test rows (Test8.png / Test9.mp4) to the Test Plan table on
# Sheet2, widens column B slightly, and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# --- 1. Stamp the new cells with the correct formatting (wrap+vcenter style
#        for the "long text" columns B/C/D, plain vcenter for E/F) by
#        copying format from existing same-styled cells. This reuses the
#        existing cellXfs entries (s="2"/s="3") instead of creating new ones.
$ws.Range("C3").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("C10").PasteSpecial(-4122)
$ws.Range("D10").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("B10").PasteSpecial(-4122)

$ws.Range("E2").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E10").PasteSpecial(-4122)

$ws.Range("F2").Copy()
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("F10").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- 2. Enter the cell values in the same order the original author typed
#        them, so new shared-string entries land at the same indices as the
#        authored workbook.
$ws.Range("B9").Value = "Buttons showing `ncorect answer"
$ws.Range("F9").Value = "Test8.png"
$ws.Range("C10").Value = "The continue buttton shoiuldn't do `nanything untill an answer is chosen. Then it will allow the user to go to the next question."
$ws.Range("C9").Value = "The 3 buttons should go red/ green`n depending on which is correct"
$ws.Range("D9").Value = "The 3 buttons go red/ green`n depending on which is correct"
$ws.Range("B10").Value = "Continue button working"
$ws.Range("F10").Value = "Test9.mp4"
$ws.Range("D10").Value = "The continue buttton doesn't`nanything."

# E9/E10 reuse the existing "Pass"/"Fail" shared strings already in the sheet.
$ws.Range("E9").Value = "Pass"
$ws.Range("E10").Value = "Fail"

# --- 3. Row heights for the two new rows (two wrapped lines / four wrapped
#        lines tall, matching the other multi-line rows on this sheet).
$ws.Rows.Item(9).RowHeight = 30
$ws.Rows.Item(10).RowHeight = 60

# --- 4. Column B needs to widen slightly to fit the new "Continue button
#        working" / "Buttons showing corect answer" text.
$ws.Columns.Item(2).ColumnWidth = 22.5

# --- 5. Move the selection/active cell (as last seen on save).
$ws.Activate()
$ws.Range("G3").Select()
